$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.512.98"
$ws.Range("E2").Value = "  -2.20%  "
$ws.Range("D3").Value = "1.995.22"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'234.05"
$ws.Range("E5").Value = "  -9.36%  "
$ws.Range("E6").Value = "  -2.41%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'54.82"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("E9").Value = "  -3.90%  "
$ws.Range("D10").Value = "'58.17"
$ws.Range("E10").Value = "  +2.61%  "
$ws.Range("E11").Value = "  -2.90%  "
$ws.Range("D12").Value = "'0.0989"
$ws.Range("E12").Value = "  -2.70%  "
$ws.Range("D13").Value = "2.286.81"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").Value = "'14.16"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "'20.35"
$ws.Range("E15").Value = "  -3.27%  "
$ws.Range("D16").Value = "'0.758"
$ws.Range("E16").Value = "  -5.42%  "
$ws.Range("E17").Value = "  -3.06%  "
$ws.Range("D18").Value = "1.997.20"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").Value = "36.469.38"
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("D20").Value = "'67.72"
$ws.Range("E20").Value = "  -2.69%  "
$ws.Range("D21").Value = "0.0₃0804"
$ws.Range("E21").Value = "  -3.63%  "
$ws.Range("E22").Value = "  +2.64%  "
$ws.Range("D23").Value = "'221.96"
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D26").Value = "'2.40"
$ws.Range("E26").Value = "  -8.63%  "
$ws.Range("D27").Value = "'161.82"
$ws.Range("E27").Value = "  -1.79%  "
$ws.Range("D28").Value = "'8.68"
$ws.Range("E28").Value = "  -2.76%  "
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("D30").Value = "'18.85"
$ws.Range("E30").Value = "  -4.11%  "
$ws.Range("D31").Value = "'1.32"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").Value = "'0.116"
$ws.Range("E32").Value = "  -2.88%  "
$ws.Range("D33").Value = "'4.39"
$ws.Range("E33").Value = "  -5.35%  "
$ws.Range("D34").Value = "'0.0602"
$ws.Range("E34").Value = "  -6.63%  "
$ws.Range("D35").Value = "'4.24"
$ws.Range("E35").Value = "  -6.85%  "
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "'3.35"
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("E39").Value = "  -3.33%  "
$ws.Range("D40").Value = "'5.59"
$ws.Range("E40").Value = "  +6.65%  "
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("D42").Value = "'0.0950"
$ws.Range("E42").Value = "  +1.97%  "
$ws.Range("D43").Value = "1.449.91"
$ws.Range("E43").Value = "  +2.97%  "
$ws.Range("E44").Value = "  -4.41%  "
$ws.Range("E45").Value = "  -8.84%  "
$ws.Range("D46").Value = "'89.01"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").Value = "'15.19"
$ws.Range("E47").Value = "  -3.35%  "
$ws.Range("D48").Value = "'0.991"
$ws.Range("E48").Value = "  -2.73%  "
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("D50").Value = "'6.84"
$ws.Range("E50").Value = "  -3.03%  "
$ws.Range("D51").Value = "'3.69"
$ws.Range("E51").Value = "  +6.59%  "
